$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of row number -> (column letter -> new text value), derived from the
# authoritative diff. Each D (Price) / E (Volume 1h) cell keeps its original
# text formatting (no real number/percentage conversion) - only the
# displayed text itself changes.
$updates = @{
    2 = @{ "D" = "298.12"; "E" = "-3.83%" }
    3 = @{ "D" = "31.70"; "E" = "-1.05%" }
    4 = @{ "D" = "5.108"; "E" = "-4.45%" }
    5 = @{ "D" = "0.07523"; "E" = "-0.66%" }
    6 = @{ "D" = "7.752"; "E" = "-0.80%" }
    7 = @{ "D" = "1.721"; "E" = "8.97%" }
    8 = @{ "D" = "3.799"; "E" = "3.43%" }
    9 = @{ "D" = "0.9305"; "E" = "2.10%" }
    10 = @{ "D" = "0.1697"; "E" = "0.45%" }
    11 = @{ "D" = "0.07503"; "E" = "-2.47%" }
    12 = @{ "D" = "0.07992"; "E" = "-1.23%" }
    13 = @{ "D" = "0.03037"; "E" = "0.51%" }
    14 = @{ "D" = "0.09891"; "E" = "0.17%" }
    15 = @{ "D" = "0.001506"; "E" = "-1.11%" }
    16 = @{ "D" = "0.006396"; "E" = "-2.00%" }
    17 = @{ "D" = "3.462"; "E" = "-1.28%" }
    18 = @{ "D" = "2.223"; "E" = "-0.68%" }
    19 = @{ "D" = "0.3274"; "E" = "0.20%" }
    20 = @{ "D" = "0.1326"; "E" = "-0.65%" }
    21 = @{ "D" = "4.559"; "E" = "8.87%" }
    22 = @{ "D" = "0.04646"; "E" = "2.11%" }
    23 = @{ "D" = "0.1558"; "E" = "-4.22%" }
    24 = @{ "D" = "0.001220"; "E" = "0.38%" }
    25 = @{ "D" = "0.004419"; "E" = "-1.57%" }
    26 = @{ "E" = "0.12%" }
    27 = @{ "E" = "6.86%" }
    39 = @{ "D" = "0.01677"; "E" = "-1.09%" }
    40 = @{ "D" = "0.04525"; "E" = "-0.67%" }
    41 = @{ "D" = "0.007060"; "E" = "-0.97%" }
    42 = @{ "D" = "0.1326"; "E" = "-2.72%" }
    43 = @{ "E" = "-8.74%" }
    44 = @{ "D" = "0.01173"; "E" = "-15.79%" }
    45 = @{ "D" = "0.00005996"; "E" = "-3.01%" }
    46 = @{ "D" = "1.918"; "E" = "1.34%" }
    47 = @{ "E" = "-0.12%" }
}

foreach ($rowKey in $updates.Keys) {
    $row = [int]$rowKey
    $cols = $updates[$rowKey]
    foreach ($col in $cols.Keys) {
        $newText = $cols[$col]
        $cell = $ws.Range("$col$row")
        # Force text storage so the numeric-looking strings ("298.12",
        # "-3.83%", etc.) are kept as literal text instead of being
        # reinterpreted as a Number/Percentage by Excel's input parser.
        $cell.NumberFormat = "@"
        $cell.Value = $newText
        # Drop the temporary "@" format again so the cell's formatting
        # stays exactly as it was before the edit (General, no style).
        $cell.ClearFormats()
    }
}
